# Auto-generated edit script for turkey_super-lig_2023-2024.xlsx
# Applies: row-level data corrections (F..V) for rows 12,13,89,90,91,92,97,98,99
# and appends three new match rows (102,103,104) with matching style for columns A & E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param($ws, $row, $data)
    foreach ($key in $data.Keys) {
        $ref = "$key$row"
        $ws.Range($ref).Value = $data[$key]
    }
}

# --- Fix mismatched match data in existing rows (home/away/odds/urls swapped) ---
Set-RowData $ws 12 @{ F="Istanbulspor AS"; G=1; H="Kayserispor"; I=1; J=1.93; K="12/08/2023 21:12"; L=2.51; M="18/08/2023 19:59"; N=3.88; O="12/08/2023 21:12"; P=3.6; Q="18/08/2023 19:59"; R=3.92; S="12/08/2023 21:12"; T=2.85; U="18/08/2023 19:59"; V="https://www.betexplorer.com/football/turkey/super-lig/istanbulspor-as-kayserispor/MLZ8zeIN/" }
Set-RowData $ws 13 @{ F="Antalyaspor"; G=1; H="Konyaspor"; I=1; J=2.32; K="14/08/2023 05:12"; L=2.05; M="18/08/2023 19:55"; N=3.57; O="14/08/2023 05:12"; P=3.6; Q="18/08/2023 19:55"; R=3.12; S="14/08/2023 05:12"; T=3.79; U="18/08/2023 19:55"; V="https://www.betexplorer.com/football/turkey/super-lig/antalyaspor-konyaspor/SzIUidmb/" }
Set-RowData $ws 89 @{ F="Trabzonspor"; G=1; H="Alanyaspor"; I=0; J=1.53; K="09/10/2023 16:12"; L=1.74; M="23/10/2023 18:58"; N=4.65; O="09/10/2023 16:12"; P=3.95; Q="23/10/2023 18:58"; R=5.9; S="09/10/2023 16:12"; T=4.95; U="23/10/2023 18:58"; V="https://www.betexplorer.com/football/turkey/super-lig/trabzonspor-alanyaspor/YePZYda9/" }
Set-RowData $ws 90 @{ F="Karagumruk"; G=2; H="Adana Demirspor"; I=0; J=3.23; K="10/10/2023 14:12"; L=3.26; M="23/10/2023 19:00"; N=3.84; O="10/10/2023 14:12"; P=3.79; Q="23/10/2023 19:00"; R=2.13; S="10/10/2023 14:12"; T=2.18; U="23/10/2023 19:00"; V="https://www.betexplorer.com/football/turkey/super-lig/f-karagumruk-adanademirspor/jXZUZxp3/" }
Set-RowData $ws 91 @{ F="Kasimpasa"; G=3; H="Istanbulspor AS"; I=1; J=1.86; K="22/10/2023 15:12"; L=1.96; M="27/10/2023 18:58"; N=3.96; O="22/10/2023 15:12"; P=3.71; Q="27/10/2023 18:59"; R=4.09; S="22/10/2023 15:12"; T=3.99; U="27/10/2023 18:58"; V="https://www.betexplorer.com/football/turkey/super-lig/kasimpasa-istanbulspor-as/YJ7MPhMe/" }
Set-RowData $ws 92 @{ F="Hatayspor"; G=1; H="Kayserispor"; I=2; J=2.26; K="22/10/2023 20:15"; L=2.49; M="27/10/2023 18:56"; N=3.55; O="22/10/2023 20:15"; P=3.42; Q="27/10/2023 18:56"; R=3.26; S="22/10/2023 20:15"; T=3.01; U="27/10/2023 18:56"; V="https://www.betexplorer.com/football/turkey/super-lig/hatayspor-kayserispor/069EREiq/" }
Set-RowData $ws 97 @{ F="Ankaragucu"; G=2; H="Samsunspor"; I=0; J=2.19; K="23/10/2023 05:42"; L=2.7; M="29/10/2023 16:59"; N=3.59; O="23/10/2023 05:42"; P=3.32; Q="29/10/2023 16:54"; R=3.29; S="23/10/2023 05:42"; T=2.82; U="29/10/2023 16:59"; V="https://www.betexplorer.com/football/turkey/super-lig/ankaragucu-samsunspor/2kdPqEDR/" }
Set-RowData $ws 98 @{ F="Antalyaspor"; G=1; H="Basaksehir"; I=0; J=2.03; K="22/10/2023 20:15"; L=2.11; M="29/10/2023 16:54"; N=3.51; O="22/10/2023 20:15"; P=3.35; Q="29/10/2023 16:54"; R=3.8; S="22/10/2023 20:15"; T=3.89; U="29/10/2023 16:54"; V="https://www.betexplorer.com/football/turkey/super-lig/antalyaspor-basaksehir/tz8GozqF/" }
Set-RowData $ws 99 @{ F="Pendikspor"; G=0; H="Fenerbahce"; I=5; J=7.87; K="23/10/2023 05:42"; L=10.56; M="29/10/2023 16:59"; N=5.68; O="23/10/2023 05:42"; P=6.22; Q="29/10/2023 16:59"; R=1.35; S="23/10/2023 05:42"; T=1.28; U="29/10/2023 16:59"; V="https://www.betexplorer.com/football/turkey/super-lig/pendikspor-fenerbahce/vc8IQY6k/" }

# --- Append new rows 102-104: copy style from the last existing data row (101) first ---
$ws.Range("A101:V101").Copy()
$ws.Range("A102:V102").PasteSpecial(-4122)
$ws.Range("A103:V103").PasteSpecial(-4122)
$ws.Range("A104:V104").PasteSpecial(-4122)

# --- Populate the new rows' values ---
Set-RowData $ws 102 @{ A=101; B="turkey"; C="super-lig"; D="2023-2024"; E=45234.47916666666; F="Istanbulspor AS"; G=2; H="Pendikspor"; I=4; J=2.08; K="30/10/2023 02:12"; L=2.67; M="04/11/2023 11:27"; N=3.77; O="30/10/2023 02:12"; P=3.53; Q="04/11/2023 11:29"; R=3.4; S="30/10/2023 02:12"; T=2.7; U="04/11/2023 11:27"; V="https://www.betexplorer.com/football/turkey/super-lig/istanbulspor-as-pendikspor/COoawCje/" }
Set-RowData $ws 103 @{ A=102; B="turkey"; C="super-lig"; D="2023-2024"; E=45234.58333333334; F="Samsunspor"; G=2; H="Hatayspor"; I=1; J=2.14; K="30/10/2023 02:12"; L=1.97; M="04/11/2023 13:55"; N=3.48; O="30/10/2023 02:12"; P=3.67; Q="04/11/2023 13:54"; R=3.51; S="30/10/2023 02:12"; T=4; U="04/11/2023 13:56"; V="https://www.betexplorer.com/football/turkey/super-lig/samsunspor-hatayspor/WQ10yjL7/" }
Set-RowData $ws 104 @{ A=103; B="turkey"; C="super-lig"; D="2023-2024"; E=45234.70833333334; F="Fenerbahce"; G=2; H="Trabzonspor"; I=3; J=1.46; K="29/10/2023 17:13"; L=1.4; M="04/11/2023 16:50"; N=4.95; O="29/10/2023 17:13"; P=5.3; Q="04/11/2023 16:59"; R=6.73; S="29/10/2023 17:13"; T=7.93; U="04/11/2023 16:59"; V="https://www.betexplorer.com/football/turkey/super-lig/fenerbahce-trabzonspor/IL6Ktocm/" }

Write-Host "Edit applied."
